$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for each row.
# Every data row (2 through 78) had its value bumped by one day: 46061 -> 46062.
for ($row = 2; $row -le 78; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
